$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" (columns A:N) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(192).Insert()

$ws1.Cells.Item(192, 1).Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws1.Cells.Item(192, 2).Value = "JUAREZ FLORES JORGE WILLIAMS"
for ($col = 3; $col -le 14; $col++) {
    $ws1.Cells.Item(192, $col).Value = 0
}

# Update the "X de 265" -> "X de 266" labels in the trailing summary row (now row 268)
for ($col = 3; $col -le 14; $col++) {
    $cell = $ws1.Cells.Item(268, $col)
    $curVal = $cell.Value()
    if ($curVal -ne $null) {
        $cell.Value = $curVal.ToString().Replace("de 265", "de 266")
    }
}

# --- Sheet 2: "VENTA MENSUAL" (columns A:G) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(192).Insert()

$ws2.Cells.Item(192, 1).Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws2.Cells.Item(192, 2).Value = "JUAREZ FLORES JORGE WILLIAMS"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(192, $col).Value = 0
}
